$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 08:08:04"
$ws1.Cells.Item(3,1).Value = "Total filas: 31"

$sheet1Rows = @(
    @("", "08:07:53", "08:10", "16_SANTA ANA",               3,  "LP1912", "30/12/2025"),
    @("", "08:07:53", "08:12", "15_ABASTO",                  5,  "LP1912", "30/12/2025"),
    @("", "08:07:53", "08:21", "26_HERNANDEZ",               14, "LP1912", "30/12/2025"),
    @("", "08:07:53", "08:22", "16_P MOR-SANTA ANA",         15, "LP1912", "30/12/2025"),
    @("", "08:07:53", "08:23", "215B_EL PATO",               16, "LP1912", "30/12/2025"),
    @("", "08:07:53", "08:27", "84_COLONIA URQUIZA-ESC 49",  20, "LP1912", "30/12/2025"),
    @("", "08:07:53", "08:33", "10_OLMOS",                   26, "LP1912", "30/12/2025"),
    @("", "08:07:53", "08:34", "16_SANTA ANA",               27, "LP1912", "30/12/2025"),
    @("", "08:07:53", "08:35", "23_HERNANDEZ",               28, "LP1912", "30/12/2025"),
    @("", "08:07:53", "08:42", "81_EL PELIGRO",              35, "LP1912", "30/12/2025"),
    @("", "08:07:53", "08:43", "14_ABASTO",                  36, "LP1912", "30/12/2025"),
    @("", "08:07:53", "08:46", "16_SANTA ANA",               39, "LP1912", "30/12/2025"),
    @("", "08:07:53", "08:53", "10_OLMOS",                   46, "LP1912", "30/12/2025"),
    @("", "08:07:53", "09:01", "215A_EL PATO",                54, "LP1912", "30/12/2025"),
    @("", "08:07:53", "09:03", "11_ETCHEVERRY",              56, "LP1912", "30/12/2025"),
    @("", "08:07:53", "09:10", "16_P MOR-SANTA ANA",         63, "LP1912", "30/12/2025"),
    @("", "08:07:53", "09:10", "23_HERNANDEZ",               63, "LP1912", "30/12/2025"),
    @("", "08:07:53", "09:13", "10_OLMOS",                   66, "LP1912", "30/12/2025"),
    @("", "08:07:53", "09:17", "27_EL RETIRO",                70, "LP1912", "30/12/2025"),
    @("", "08:07:53", "09:21", "26_HERNANDEZ",               74, "LP1912", "30/12/2025"),
    @("", "08:07:53", "09:23", "11_ETCHEVERRY",              76, "LP1912", "30/12/2025"),
    @("", "08:07:53", "09:42", "215C_EL PATO",               95, "LP1912", "30/12/2025"),
    @("", "08:07:53", "09:42", "23_HERNANDEZ",               95, "LP1912", "30/12/2025")
)

$r = 10
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 08:08:04"
$ws2.Cells.Item(3,1).Value = "Total filas: 6"

$sheet2Rows = @(
    @("", "30/12/2025", "08:07:53", "08:23", "215B_EL PATO", 16, "LP1912"),
    @("", "30/12/2025", "08:07:53", "09:01", "215A_EL PATO", 54, "LP1912"),
    @("", "30/12/2025", "08:07:53", "09:42", "215C_EL PATO", 95, "LP1912")
)

$r = 5
foreach ($row in $sheet2Rows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $ws2.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 08:08:04"
$ws3.Cells.Item(3,1).Value = "Total filas: 5"

$sheet3Rows = @(
    @("", "30/12/2025", "08:08:04", "08:39", "215A_LA PLATA", 31, "L6173"),
    @("", "30/12/2025", "08:07:59", "09:08", "215D_LA PLATA", 61, "L6203")
)

$r = 5
foreach ($row in $sheet3Rows) {
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $ws3.Cells.Item($r, 5).Value = $row[4]
    $ws3.Cells.Item($r, 6).Value = $row[5]
    $ws3.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}
